# Adapt column header formatting to respective input file names:
#   "<header>_old" -> "<header>_FV2404"
#   "<header>_new" -> "<header>_FV2410"
# Then turn the data range into a proper Excel Table (ListObject) and
# freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the header row -------------------------------------------------
# Columns A:J carry the "_old" suffixed headers, column K is "diff", and
# columns L:U carry the "_new" suffixed headers.
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $name = $cell.Value()
    if ($name -like "*_old") {
        $base = $name.Substring(0, $name.Length - 4)
        $cell.Value = $base + "_FV2404"
    } elseif ($name -like "*_new") {
        $base = $name.Substring(0, $name.Length - 4)
        $cell.Value = $base + "_FV2410"
    }
}

# --- 2) Turn the used range into a Table (ListObject) -------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$tableRange = $ws.Range("A1:U" + $lastRow)
$tbl = $ws.ListObjects.Add(1, $tableRange, [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ---------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
